# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.671.49'; E = '  -0.29%  ' }
    @{ Row = 3; D = '1.583.37'; E = '  -2.82%  ' }
    @{ Row = 4; D = $null; E = '  +0.70%  ' }
    @{ Row = 5; D = '206.80'; E = '  -1.92%  ' }
    @{ Row = 6; D = '0.507'; E = '  -2.21%  ' }
    @{ Row = 7; D = $null; E = '  +0.79%  ' }
    @{ Row = 8; D = $null; E = '  -4.52%  ' }
    @{ Row = 9; D = $null; E = '  -1.53%  ' }
    @{ Row = 10; D = '0.0591'; E = '  -2.95%  ' }
    @{ Row = 11; D = $null; E = '  -1.32%  ' }
    @{ Row = 12; D = '1.808.56'; E = '  -3.05%  ' }
    @{ Row = 13; D = '1.575.67'; E = '  -3.55%  ' }
    @{ Row = 14; D = $null; E = '  -3.55%  ' }
    @{ Row = 15; D = $null; E = '  -5.41%  ' }
    @{ Row = 16; D = '27.662.17'; E = '  -0.74%  ' }
    @{ Row = 17; D = '63.26'; E = '  -2.83%  ' }
    @{ Row = 18; D = '219.53'; E = '  -4.17%  ' }
    @{ Row = 19; D = '0.0₃0693'; E = '  -3.49%  ' }
    @{ Row = 20; D = '7.33'; E = '  -5.59%  ' }
    @{ Row = 21; D = $null; E = '  +0.55%  ' }
    @{ Row = 22; D = $null; E = '  -4.71%  ' }
    @{ Row = 23; D = '9.51'; E = '  -5.82%  ' }
    @{ Row = 24; D = $null; E = '  -5.58%  ' }
    @{ Row = 25; D = '154.14'; E = '  -0.69%  ' }
    @{ Row = 26; D = '6.79'; E = '  -2.35%  ' }
    @{ Row = 27; D = $null; E = '  +0.61%  ' }
    @{ Row = 28; D = $null; E = '  -2.48%  ' }
    @{ Row = 29; D = $null; E = '  -3.96%  ' }
    @{ Row = 30; D = $null; E = '  -2.05%  ' }
    @{ Row = 31; D = $null; E = '  -3.18%  ' }
    @{ Row = 32; D = $null; E = '  -5.10%  ' }
    @{ Row = 33; D = '1.384.18'; E = '  -1.03%  ' }
    @{ Row = 34; D = '2.92'; E = '  -5.51%  ' }
    @{ Row = 35; D = $null; E = '  -5.24%  ' }
    @{ Row = 36; D = $null; E = '  -4.65%  ' }
    @{ Row = 37; D = $null; E = '  -0.67%  ' }
    @{ Row = 38; D = $null; E = '  -3.55%  ' }
    @{ Row = 39; D = '0.539'; E = '  -3.27%  ' }
    @{ Row = 40; D = '0.821'; E = '  -3.21%  ' }
    @{ Row = 41; D = $null; E = '  +0.59%  ' }
    @{ Row = 42; D = '0.976'; E = '  -3.78%  ' }
    @{ Row = 43; D = '63.47'; E = '  -3.65%  ' }
    @{ Row = 44; D = $null; E = '  -0.46%  ' }
    @{ Row = 45; D = '1.75'; E = '  -3.95%  ' }
    @{ Row = 46; D = '5.22'; E = '  -3.99%  ' }
    @{ Row = 47; D = '1.719.96'; E = '  -3.10%  ' }
    @{ Row = 48; D = '88.14'; E = '  -0.13%  ' }
    @{ Row = 49; D = $null; E = '  +5.32%  ' }
    @{ Row = 50; D = '0.0972'; E = '  -4.70%  ' }
    @{ Row = 51; D = $null; E = '  -1.17%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)   # column D
        # Force plain-text storage so values that look numeric (e.g. "206.80")
        # keep their exact original digits instead of being parsed as numbers.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $eCell = $ws.Cells.Item($u.Row, 5)   # column E
        $eCell.Value = $u.E
    }
}

